$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arrB = New-Object 'object[,]' 24,1
$arrB[0,0] = 0.9752684689550506
$arrB[1,0] = 0.8903146079812245
$arrB[2,0] = 0.8383503803046324
$arrB[3,0] = 0.8172253936473055
$arrB[4,0] = 0.813720710537666
$arrB[5,0] = 0.8380652739348591
$arrB[6,0] = 0.9459359731671384
$arrB[7,0] = 1.159000976921334
$arrB[8,0] = 1.31643747657813
$arrB[9,0] = 1.38824777267132
$arrB[10,0] = 1.415467106020287
$arrB[11,0] = 1.409603783902298
$arrB[12,0] = 1.390486602873239
$arrB[13,0] = 1.378780170457219
$arrB[14,0] = 1.311748256403689
$arrB[15,0] = 1.270674560187388
$arrB[16,0] = 1.247068167924624
$arrB[17,0] = 1.239078597143532
$arrB[18,0] = 1.275045059683123
$arrB[19,0] = 1.396101077442495
$arrB[20,0] = 1.475371102920349
$arrB[21,0] = 1.433049610649164
$arrB[22,0] = 1.273069133646686
$arrB[23,0] = 1.101201047570214
$ws.Range("B2:B25").Value = $arrB

$arrC = New-Object 'object[,]' 24,1
$arrC[0,0] = 0.3144284655096214
$arrC[1,0] = 0.3048098054362356
$arrC[2,0] = 0.2988588799797185
$arrC[3,0] = 0.2964226609837226
$arrC[4,0] = 0.2960174586737168
$arrC[5,0] = 0.2988260692945062
$arrC[6,0] = 0.3111213871380443
$arrC[7,0] = 0.3348694634296407
$arrC[8,0] = 0.3520897647377126
$arrC[9,0] = 0.3598731472398526
$arrC[10,0] = 0.3628131539281583
$arrC[11,0] = 0.3621803023158918
$arrC[12,0] = 0.3601151724663509
$arrC[13,0] = 0.3588492531960412
$arrC[14,0] = 0.3515800773149067
$arrC[15,0] = 0.3471076900233356
$arrC[16,0] = 0.3445305744290863
$arrC[17,0] = 0.3436572024542102
$arrC[18,0] = 0.3475842726818712
$arrC[19,0] = 0.3607219531134547
$arrC[20,0] = 0.3692650258070671
$arrC[21,0] = 0.364709432373985
$arrC[22,0] = 0.3473688279672444
$arrC[23,0] = 0.328484475510038
$ws.Range("C2:C25").Value = $arrC

$arrD = New-Object 'object[,]' 24,1
$arrD[0,0] = 0.02482583923178794
$arrD[1,0] = 0.02364819138880847
$arrD[2,0] = 0.02291906668793686
$arrD[3,0] = 0.0226204400608232
$arrD[4,0] = 0.02257076311077455
$arrD[5,0] = 0.02291504536507816
$arrD[6,0] = 0.02442105085827251
$arrD[7,0] = 0.02732572685741275
$arrD[8,0] = 0.0294294982446317
$arrD[9,0] = 0.03037984995191323
$arrD[10,0] = 0.03073874974306534
$arrD[11,0] = 0.03066149804794094
$arrD[12,0] = 0.03040939656172981
$arrD[13,0] = 0.03025484921101906
$arrD[14,0] = 0.02936725475735358
$arrD[15,0] = 0.02882102354342919
$arrD[16,0] = 0.02850621966682354
$arrD[17,0] = 0.0283995255377647
$arrD[18,0] = 0.02887923573199913
$arrD[19,0] = 0.03048347157044162
$arrD[20,0] = 0.0315262242068215
$arrD[21,0] = 0.03097021621596951
$arrD[22,0] = 0.02885292039287179
$arrD[23,0] = 0.02654520684647821
$ws.Range("D2:D25").Value = $arrD

$arrF = New-Object 'object[,]' 24,1
$arrF[0,0] = 0.7747881010662638
$arrF[1,0] = 0.7683249196520308
$arrF[2,0] = 0.764887090079732
$arrF[3,0] = 0.7636194630226427
$arrF[4,0] = 0.7634170234460527
$arrF[5,0] = 0.7648694547651118
$arrF[6,0] = 0.7724494048240231
$arrF[7,0] = 0.7915303816012624
$arrF[8,0] = 0.8081319906748377
$arrF[9,0] = 0.8162481181249177
$arrF[10,0] = 0.8194027483387885
$arrF[11,0] = 0.8187197275378253
$arrF[12,0] = 0.8165060227605636
$arrF[13,0] = 0.8151606469139807
$arrF[14,0] = 0.807612942392538
$arrF[15,0] = 0.8031272100857905
$arrF[16,0] = 0.8006002115151745
$arrF[17,0] = 0.7997537250854805
$arrF[18,0] = 0.8035992298306809
$arrF[19,0] = 0.8171540359525409
$arrF[20,0] = 0.8264864091584911
$arrF[21,0] = 0.8214621757054488
$arrF[22,0] = 0.80338566802088
$arrF[23,0] = 0.7859157796321767
$ws.Range("F2:F25").Value = $arrF

$arrG = New-Object 'object[,]' 24,1
$arrG[0,0] = 0.002426413660574234
$arrG[1,0] = 0.002429291666281998
$arrG[2,0] = 0.002431154550266353
$arrG[3,0] = 0.002431937846937425
$arrG[4,0] = 0.002432069373864594
$arrG[5,0] = 0.0024311650164509
$arrG[6,0] = 0.002427386162088477
$arrG[7,0] = 0.002420732487442757
$arrG[8,0] = 0.002416300685649583
$arrG[9,0] = 0.002414382719548659
$arrG[10,0] = 0.002413670465668276
$arrG[11,0] = 0.002413823239018478
$arrG[12,0] = 0.002414323841150414
$arrG[13,0] = 0.002414632300641887
$arrG[14,0] = 0.002416427997163303
$arrG[15,0] = 0.002417554673284417
$arrG[16,0] = 0.00241821194328592
$arrG[17,0] = 0.002418436072071076
$arrG[18,0] = 0.00241743378113608
$arrG[19,0] = 0.002414176421762803
$arrG[20,0] = 0.002412129343998635
$arrG[21,0] = 0.002413214445327435
$arrG[22,0] = 0.002417488406655779
$arrG[23,0] = 0.002422451957877659
$ws.Range("G2:G25").Value = $arrG

$arrL = New-Object 'object[,]' 24,1
$arrL[0,0] = 0.2699937045276712
$arrL[1,0] = 0.2680661493647136
$arrL[2,0] = 0.2670216154850706
$arrL[3,0] = 0.2666309506830018
$arrL[4,0] = 0.266568195669258
$arrL[5,0] = 0.2670162051111049
$arrL[6,0] = 0.2693002584326791
$arrL[7,0] = 0.2748809283151985
$arrL[8,0] = 0.2796518848417691
$arrL[9,0] = 0.2819678800110523
$arrL[10,0] = 0.2828658072999275
$arrL[11,0] = 0.2826714931794214
$arrL[12,0] = 0.2820413341644894
$arrL[13,0] = 0.2816580656271412
$arrL[14,0] = 0.2795034571476123
$arrL[15,0] = 0.2782189568330864
$arrL[16,0] = 0.2774938588810869
$arrL[17,0] = 0.2772507094320673
$arrL[18,0] = 0.2783542751377297
$arrL[19,0] = 0.2822258599331917
$arrL[20,0] = 0.2848780282515122
$arrL[21,0] = 0.2834513774452745
$arrL[22,0] = 0.2782930560581889
$arrL[23,0] = 0.2732533486240172
$ws.Range("L2:L25").Value = $arrL

$arrM = New-Object 'object[,]' 24,1
$arrM[0,0] = 0.2426022336133826
$arrM[1,0] = 0.2293245913893287
$arrM[2,0] = 0.2212654132829854
$arrM[3,0] = 0.2180048774441445
$arrM[4,0] = 0.2174649010964913
$arrM[5,0] = 0.2212213445922018
$arrM[6,0] = 0.2380048265235075
$arrM[7,0] = 0.2716521723259717
$arrM[8,0] = 0.296816041549981
$arrM[9,0] = 0.3083590973021373
$arrM[10,0] = 0.3127438058863277
$arrM[11,0] = 0.3117988783364396
$arrM[12,0] = 0.3087195583291589
$arrM[13,0] = 0.3068351520718338
$arrM[14,0] = 0.2960635894091013
$arrM[15,0] = 0.2894800119561367
$arrM[16,0] = 0.2857023485304566
$arrM[17,0] = 0.2844248541987824
$arrM[18,0] = 0.290179911138587
$arrM[19,0] = 0.3096236615215275
$arrM[20,0] = 0.3224104848702183
$arrM[21,0] = 0.3155787298721719
$arrM[22,0] = 0.2898634638943705
$arrM[23,0] = 0.2624714962310435
$ws.Range("M2:M25").Value = $arrM

$arrN = New-Object 'object[,]' 24,1
$arrN[0,0] = 1.359822563877721
$arrN[1,0] = 1.374288397553542
$arrN[2,0] = 1.383704108734388
$arrN[3,0] = 1.387675299132745
$arrN[4,0] = 1.388342819160151
$arrN[5,0] = 1.383757122144274
$arrN[6,0] = 1.364699644013445
$arrN[7,0] = 1.331560704272803
$arrN[8,0] = 1.309791617105091
$arrN[9,0] = 1.300447685975541
$arrN[10,0] = 1.296989754401949
$arrN[11,0] = 1.297730904824078
$arrN[12,0] = 1.30016158804731
$arrN[13,0] = 1.301660924970356
$arrN[14,0] = 1.31041351535378
$arrN[15,0] = 1.315926128667407
$arrN[16,0] = 1.319149445170417
$arrN[17,0] = 1.320249838992474
$arrN[18,0] = 1.315333856663926
$arrN[19,0] = 1.299445454771337
$arrN[20,0] = 1.289530191057601
$arrN[21,0] = 1.294779251833376
$arrN[22,0] = 1.315601454499152
$arrN[23,0] = 1.340072777324981
$ws.Range("N2:N25").Value = $arrN

$arrO = New-Object 'object[,]' 24,1
$arrO[0,0] = 2.660753901150883
$arrO[1,0] = 2.653760480629671
$arrO[2,0] = 2.651247702501365
$arrO[3,0] = 2.650671200377644
$arrO[4,0] = 2.650602486669413
$arrO[5,0] = 2.651238116349361
$arrO[6,0] = 2.657972613782107
$arrO[7,0] = 2.685336084649009
$arrO[8,0] = 2.714110925082167
$arrO[9,0] = 2.72909326382387
$arrO[10,0] = 2.735039427249234
$arrO[11,0] = 2.733746678689897
$arrO[12,0] = 2.72957699035976
$arrO[13,0] = 2.727058462281263
$arrO[14,0] = 2.713169922114076
$arrO[15,0] = 2.70513485348701
$arrO[16,0] = 2.700691395561137
$arrO[17,0] = 2.69921748975608
$arrO[18,0] = 2.705971763584671
$arrO[19,0] = 2.730794324271329
$arrO[20,0] = 2.748606919051298
$arrO[21,0] = 2.738954375786363
$arrO[22,0] = 2.705592848303695
$arrO[23,0] = 2.676413748109439
$ws.Range("O2:O25").Value = $arrO

Write-Host "applied changes"